$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells are not auto-converted to numbers/dates by Excel
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.628.97"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "2.538.10"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "594.10"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "178.01"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "2.537.21"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +5.46%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").Value = "5.01"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "0.338"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "2.983.77"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "26.15"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "68.533.22"
$ws.Range("E16").Value = "  +2.19%  "
$ws.Range("D17").Value = "0.0000171"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "2.537.86"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").Value = "11.13"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "353.17"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "4.22"
$ws.Range("E22").Value = "  +4.68%  "
$ws.Range("D23").Value = "71.34"
$ws.Range("E23").Value = "  +2.80%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "4.23"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "2.639.10"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").Value = "516.19"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("D31").Value = "0.0₃0903"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "7.84"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("E33").Value = "  +2.63%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").Value = "163.05"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "18.45"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "18.70"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.76"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "4.86"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.328"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").Value = "152.62"
$ws.Range("E46").Value = "  +7.12%  "
$ws.Range("D47").Value = "3.58"
$ws.Range("E47").Value = "  +2.73%  "
$ws.Range("D48").Value = "0.522"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "1.62"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0743"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.584"
$ws.Range("E51").Value = "  +0.36%  "
